# Refresh the cryptos snapshot (price + 1h volume columns, plus three
# coin rows whose ranking moved) to match the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new literal text. Columns D (Price)
# and E (Volume(1h)) must stay plain text -- many of the values look like
# numbers (e.g. "1.00", "422.21") or would otherwise lose formatting (the
# thousands-dot style prices, the padded "  +x.xx%  " volumes) if Excel
# auto-coerced them, so NumberFormat is forced to Text ("@") first.
$updates = @(
    @{ Cell = "D2"; Value = "66.183.58" }
    @{ Cell = "E2"; Value = "  -1.72%  " }
    @{ Cell = "D3"; Value = "3.830.24" }
    @{ Cell = "E3"; Value = "  +1.90%  " }
    @{ Cell = "D4"; Value = "0.998" }
    @{ Cell = "E4"; Value = "  -0.40%  " }
    @{ Cell = "D5"; Value = "422.21" }
    @{ Cell = "E5"; Value = "  +0.38%  " }
    @{ Cell = "D6"; Value = "129.17" }
    @{ Cell = "E6"; Value = "  -2.13%  " }
    @{ Cell = "D7"; Value = "3.911.07" }
    @{ Cell = "E7"; Value = "  +4.26%  " }
    @{ Cell = "D8"; Value = "0.603" }
    @{ Cell = "E8"; Value = "  -7.28%  " }
    @{ Cell = "E9"; Value = "  +0.07%  " }
    @{ Cell = "D10"; Value = "0.717" }
    @{ Cell = "E10"; Value = "  -7.47%  " }
    @{ Cell = "D11"; Value = "0.163" }
    @{ Cell = "E11"; Value = "  -12.44%  " }
    @{ Cell = "D12"; Value = "0.0000344" }
    @{ Cell = "E12"; Value = "  -18.36%  " }
    @{ Cell = "D13"; Value = "40.29" }
    @{ Cell = "E13"; Value = "  -6.15%  " }
    @{ Cell = "D14"; Value = "4.425.67" }
    @{ Cell = "E14"; Value = "  +1.65%  " }
    @{ Cell = "B15"; Value = "Polkadot" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" }
    @{ Cell = "D15"; Value = "9.96" }
    @{ Cell = "E15"; Value = "  -4.66%  " }
    @{ Cell = "B16"; Value = "Uniswap" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" }
    @{ Cell = "D16"; Value = "15.88" }
    @{ Cell = "E16"; Value = "  +20.70%  " }
    @{ Cell = "B17"; Value = "WrappedEther" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell = "D17"; Value = "3.830.64" }
    @{ Cell = "E17"; Value = "  +2.28%  " }
    @{ Cell = "B18"; Value = "TRON" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx" }
    @{ Cell = "D18"; Value = "0.137" }
    @{ Cell = "E18"; Value = "  -1.63%  " }
    @{ Cell = "D19"; Value = "19.43" }
    @{ Cell = "E19"; Value = "  -5.95%  " }
    @{ Cell = "D20"; Value = "66.319.03" }
    @{ Cell = "E20"; Value = "  -1.53%  " }
    @{ Cell = "E21"; Value = "  -7.42%  " }
    @{ Cell = "D22"; Value = "400.79" }
    @{ Cell = "E22"; Value = "  -10.82%  " }
    @{ Cell = "D23"; Value = "14.32" }
    @{ Cell = "E23"; Value = "  -10.14%  " }
    @{ Cell = "D24"; Value = "83.76" }
    @{ Cell = "E24"; Value = "  -6.31%  " }
    @{ Cell = "D25"; Value = "2.99" }
    @{ Cell = "E25"; Value = "  -3.34%  " }
    @{ Cell = "D26"; Value = "36.82" }
    @{ Cell = "E26"; Value = "  -5.03%  " }
    @{ Cell = "E27"; Value = "  +11.75%  " }
    @{ Cell = "D28"; Value = "3.21" }
    @{ Cell = "E28"; Value = "  -3.93%  " }
    @{ Cell = "D29"; Value = "9.37" }
    @{ Cell = "E29"; Value = "  -7.91%  " }
    @{ Cell = "D30"; Value = "692.00" }
    @{ Cell = "E30"; Value = "  +1.78%  " }
    @{ Cell = "D31"; Value = "2.75" }
    @{ Cell = "E31"; Value = "  +1.64%  " }
    @{ Cell = "D32"; Value = "0.120" }
    @{ Cell = "E32"; Value = "  -4.62%  " }
    @{ Cell = "D33"; Value = "12.25" }
    @{ Cell = "E33"; Value = "  -3.76%  " }
    @{ Cell = "D34"; Value = "7.30" }
    @{ Cell = "E34"; Value = "  +0.23%  " }
    @{ Cell = "E35"; Value = "  -9.63%  " }
    @{ Cell = "D36"; Value = "37.89" }
    @{ Cell = "E36"; Value = "  -9.69%  " }
    @{ Cell = "D37"; Value = "1.00" }
    @{ Cell = "E37"; Value = "  +0.00%  " }
    @{ Cell = "D38"; Value = "54.67" }
    @{ Cell = "E38"; Value = "  -4.28%  " }
    @{ Cell = "D39"; Value = "0.0₃0761" }
    @{ Cell = "E39"; Value = "  -1.89%  " }
    @{ Cell = "D40"; Value = "0.0451" }
    @{ Cell = "E40"; Value = "  -8.55%  " }
    @{ Cell = "D41"; Value = "2.94" }
    @{ Cell = "E41"; Value = "  -0.09%  " }
    @{ Cell = "D42"; Value = "0.996" }
    @{ Cell = "E42"; Value = "  -0.07%  " }
    @{ Cell = "E43"; Value = "  -10.38%  " }
    @{ Cell = "D44"; Value = "4.46" }
    @{ Cell = "E44"; Value = "  +2.45%  " }
    @{ Cell = "D45"; Value = "143.99" }
    @{ Cell = "E45"; Value = "  -3.05%  " }
    @{ Cell = "B46"; Value = "LidoDAOToken" }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo" }
    @{ Cell = "D46"; Value = "3.25" }
    @{ Cell = "E46"; Value = "  -5.35%  " }
    @{ Cell = "B47"; Value = "EnergySwap" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D47"; Value = "26.27" }
    @{ Cell = "E47"; Value = "  -7.38%  " }
    @{ Cell = "B48"; Value = "ApeXProtocol" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex" }
    @{ Cell = "D48"; Value = "3.07" }
    @{ Cell = "E48"; Value = "  -3.08%  " }
    @{ Cell = "D49"; Value = "2.04" }
    @{ Cell = "E49"; Value = "  -5.03%  " }
    @{ Cell = "D50"; Value = "2.52" }
    @{ Cell = "E50"; Value = "  -4.64%  " }
    @{ Cell = "D51"; Value = "2.74" }
    @{ Cell = "E51"; Value = "  -5.32%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell -match "^[DE]") {
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
